$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 695.81354
$ws.Range("J17").Value = 695.81354
$ws.Range("L17").Value = 2087.44062
$ws.Range("N17").Value = -2423.44062

$ws.Range("H62").Value = 3688.75
$ws.Range("I62").Value = 3252
$ws.Range("K62").Value = 3252
$ws.Range("M62").Value = -2628

$ws.Range("H65").Value = 3688.75
$ws.Range("I65").Value = 3252
$ws.Range("K65").Value = 16260
$ws.Range("M65").Value = -13140

$ws.Range("H96").Value = 353.25
$ws.Range("I96").Value = 376.36365
$ws.Range("J96").Value = 99
$ws.Range("K96").Value = 1129.09095
$ws.Range("L96").Value = 297
$ws.Range("M96").Value = 243.90905
$ws.Range("N96").Value = -3043

$ws.Range("H112").Value = 3388.375
$ws.Range("I112").Value = 1989
$ws.Range("K112").Value = 5967
$ws.Range("M112").Value = -4859

$ws.Range("H138").Value = 32260112
$ws.Range("J138").Value = 62502028
$ws.Range("L138").Value = 187506084
$ws.Range("N138").Value = -187516364

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4740.6265
$ws.Range("I32").Value = 4717.5776
$ws.Range("K32").Value = 4717.5776
$ws.Range("M32").Value = -4430.5776

$ws.Range("H44").Value = 31000
$ws.Range("J44").Value = 31000
$ws.Range("L44").Value = 31000
$ws.Range("N44").Value = -31976

$ws.Range("H55").Value = 25737.5
$ws.Range("J55").Value = 36500
$ws.Range("L55").Value = 36500
$ws.Range("N55").Value = -37130

$ws.Range("H80").Value = 41950
$ws.Range("J80").Value = 41950
$ws.Range("L80").Value = 41950
$ws.Range("N80").Value = -43946

$ws.Range("H83").Value = 41950
$ws.Range("J83").Value = 41950
$ws.Range("L83").Value = 125850
$ws.Range("N83").Value = -135834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 3952.9167
$ws.Range("I94").Value = 3952.9167
$ws.Range("K94").Value = 3952.9167
$ws.Range("M94").Value = -3501.9167

$ws.Range("H107").Value = 1304.6
$ws.Range("I107").Value = 1267.0588
$ws.Range("K107").Value = 1267.0588
$ws.Range("M107").Value = 652.9412

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 22699.5
$ws.Range("J41").Value = 21200
$ws.Range("L41").Value = 21200
$ws.Range("N41").Value = -22056

$ws.Range("H59").Value = 40747.25
$ws.Range("J59").Value = 40996.332
$ws.Range("L59").Value = 40996.332
$ws.Range("N59").Value = -43286.332

$ws.Range("H60").Value = 24500
$ws.Range("J60").Value = 24500
$ws.Range("L60").Value = 24500
$ws.Range("N60").Value = -25522

$ws.Range("H74").Value = 37457
$ws.Range("J74").Value = 37457
$ws.Range("L74").Value = 37457
$ws.Range("N74").Value = -39205

$ws.Range("H77").Value = 37457
$ws.Range("J77").Value = 37457
$ws.Range("L77").Value = 112371
$ws.Range("N77").Value = -121107

$ws.Range("H132").Value = 3378.9285
$ws.Range("I132").Value = 3485
$ws.Range("K132").Value = 10455
$ws.Range("M132").Value = -7925

$ws.Range("H134").Value = 4571.457
$ws.Range("I134").Value = 3537.0908
$ws.Range("J134").Value = 6321.923
$ws.Range("K134").Value = 10611.2724
$ws.Range("L134").Value = 18965.769
$ws.Range("M134").Value = -8076.2724
$ws.Range("N134").Value = -24035.769

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 1260.6
$ws.Range("J12").Value = 1563.3636
$ws.Range("L12").Value = 4690.0908
$ws.Range("N12").Value = -5036.0908

$ws.Range("H68").Value = 3404.818
$ws.Range("I68").Value = 678.1429000000001
$ws.Range("J68").Value = 8176.5
$ws.Range("K68").Value = 2034.4287
$ws.Range("L68").Value = 24529.5
$ws.Range("M68").Value = -1223.4287
$ws.Range("N68").Value = -26151.5

$ws.Range("H71").Value = 3404.818
$ws.Range("I71").Value = 678.1429000000001
$ws.Range("J71").Value = 8176.5
$ws.Range("K71").Value = 6103.2861
$ws.Range("L71").Value = 73588.5
$ws.Range("M71").Value = -2047.2861
$ws.Range("N71").Value = -81700.5

$ws.Range("H103").Value = 165
$ws.Range("I103").Value = 165
$ws.Range("K103").Value = 495
$ws.Range("M103").Value = 384

$ws.Range("H122").Value = 731.25
$ws.Range("J122").Value = 811.6667
$ws.Range("L122").Value = 7305.0003
$ws.Range("N122").Value = -12205.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1669684
$ws.Range("I3").Value = 5000500
$ws.Range("K3").Value = 5000500
$ws.Range("M3").Value = -5000384

$ws.Range("H80").Value = 3627
$ws.Range("I80").Value = 3679.25
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 3679.25
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -2681.25
$ws.Range("N80").Value = -4996

$ws.Range("H83").Value = 3627
$ws.Range("I83").Value = 3679.25
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 18396.25
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -13404.25
$ws.Range("N83").Value = -24984

$ws.Range("H107").Value = 656.1111
$ws.Range("I107").Value = 809.7
$ws.Range("J107").Value = 464.125
$ws.Range("K107").Value = 809.7
$ws.Range("L107").Value = 464.125
$ws.Range("M107").Value = 1110.3
$ws.Range("N107").Value = -4304.125

$ws.Range("H113").Value = 3512.2856
$ws.Range("I113").Value = 4166.6665
$ws.Range("J113").Value = 3021.5
$ws.Range("K113").Value = 4166.6665
$ws.Range("L113").Value = 3021.5
$ws.Range("M113").Value = -1996.6665
$ws.Range("N113").Value = -7361.5

$ws.Range("H123").Value = 42068
$ws.Range("J123").Value = 42068
$ws.Range("L123").Value = 42068
$ws.Range("N123").Value = -46968

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2601.8
$ws.Range("I7").Value = 2501
$ws.Range("J7").Value = 3005
$ws.Range("K7").Value = 2501
$ws.Range("L7").Value = 3005
$ws.Range("M7").Value = -2389
$ws.Range("N7").Value = -3229

$ws.Range("H42").Value = 11984
$ws.Range("I42").Value = 8999
$ws.Range("J42").Value = 14969
$ws.Range("K42").Value = 8999
$ws.Range("L42").Value = 14969
$ws.Range("M42").Value = -8436
$ws.Range("N42").Value = -16095

$ws.Range("H49").Value = 11984
$ws.Range("I49").Value = 8999
$ws.Range("J49").Value = 14969
$ws.Range("K49").Value = 8999
$ws.Range("L49").Value = 14969
$ws.Range("M49").Value = -8852
$ws.Range("N49").Value = -15263

$ws.Range("H61").Value = 51236.5
$ws.Range("I61").Value = 56887.61
$ws.Range("J61").Value = 376.5
$ws.Range("K61").Value = 56887.61
$ws.Range("L61").Value = 376.5
$ws.Range("M61").Value = -56685.61
$ws.Range("N61").Value = -780.5

$ws.Range("H68").Value = 1841.5
$ws.Range("I68").Value = 1837.5
$ws.Range("J68").Value = 1849.5
$ws.Range("K68").Value = 1837.5
$ws.Range("L68").Value = 1849.5
$ws.Range("M68").Value = -1088.5
$ws.Range("N68").Value = -3347.5

$ws.Range("H71").Value = 1841.5
$ws.Range("I71").Value = 1837.5
$ws.Range("J71").Value = 1849.5
$ws.Range("K71").Value = 9187.5
$ws.Range("L71").Value = 9247.5
$ws.Range("M71").Value = -5443.5
$ws.Range("N71").Value = -16735.5

$ws.Range("H93").Value = 2033
$ws.Range("I93").Value = 2013.25
$ws.Range("K93").Value = 2013.25
$ws.Range("M93").Value = -765.25

$ws.Range("H113").Value = 51236.5
$ws.Range("I113").Value = 56887.61
$ws.Range("J113").Value = 376.5
$ws.Range("K113").Value = 56887.61
$ws.Range("L113").Value = 376.5
$ws.Range("M113").Value = -54717.61
$ws.Range("N113").Value = -4716.5

$ws.Range("H126").Value = 2601.8
$ws.Range("I126").Value = 2501
$ws.Range("J126").Value = 3005
$ws.Range("K126").Value = 7503
$ws.Range("L126").Value = 9015
$ws.Range("M126").Value = -5033
$ws.Range("N126").Value = -13955

$ws.Range("H132").Value = 7773.1904
$ws.Range("I132").Value = 8080.972
$ws.Range("K132").Value = 24242.916
$ws.Range("M132").Value = -21712.916

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 10933.125
$ws.Range("J5").Value = 12280.714
$ws.Range("L5").Value = 12280.714
$ws.Range("N5").Value = -12504.714

$ws.Range("H25").Value = 8030
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = ""

$ws.Range("H54").Value = 29500
$ws.Range("J54").Value = 29500
$ws.Range("L54").Value = 29500
$ws.Range("N54").Value = -30540

$ws.Range("H70").Value = 40492
$ws.Range("J70").Value = 40492
$ws.Range("L70").Value = 40492
$ws.Range("N70").Value = -41122

$ws.Range("H73").Value = 40492
$ws.Range("J73").Value = 40492
$ws.Range("L73").Value = 40492
$ws.Range("N73").Value = -42676

$ws.Range("H81").Value = 5949.6665
$ws.Range("I81").Value = 5139.2
$ws.Range("K81").Value = 10278.4
$ws.Range("M81").Value = -9217.4

$ws.Range("H84").Value = 5949.6665
$ws.Range("I84").Value = 5139.2
$ws.Range("K84").Value = 51392
$ws.Range("M84").Value = -46088

$ws.Range("H126").Value = 2005.1428
$ws.Range("I126").Value = 2005.1428
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6015.428400000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = -3545.428400000001

$ws.Range("H136").Value = 3617.725
$ws.Range("I136").Value = 2948.8286
$ws.Range("K136").Value = 8846.485799999999
$ws.Range("M136").Value = -6296.485799999999
